$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 333
$ws.Range("I4").Value = 183.42857
$ws.Range("K4").Value = 183.42857
$ws.Range("M4").Value = -69.42857000000001

$ws.Range("H100").Value = 7543.227
$ws.Range("J100").Value = 10492.066
$ws.Range("L100").Value = 10492.066
$ws.Range("N100").Value = -11574.066

$ws.Range("H112").Value = 772035.0600000001
$ws.Range("J112").Value = 836304.7
$ws.Range("L112").Value = 2508914.1
$ws.Range("N112").Value = -2511130.1

$ws.Range("H135").Value = 2853.4688
$ws.Range("I135").Value = 1259.1364
$ws.Range("K135").Value = 11332.2276
$ws.Range("M135").Value = -8797.2276

$ws.Range("H136").Value = 55126.605
$ws.Range("J136").Value = 50909.09
$ws.Range("L136").Value = 50909.09
$ws.Range("N136").Value = -61109.09

$ws.Range("H137").Value = 4975.9644
$ws.Range("I137").Value = 4611.048
$ws.Range("J137").Value = 6070.7144
$ws.Range("K137").Value = 13833.144
$ws.Range("L137").Value = 18212.1432
$ws.Range("M137").Value = -11283.144
$ws.Range("N137").Value = -23312.1432

$ws.Range("H138").Value = 4139.91
$ws.Range("I138").Value = 1314.6923
$ws.Range("J138").Value = 4562.069
$ws.Range("K138").Value = 3944.0769
$ws.Range("L138").Value = 13686.207
$ws.Range("M138").Value = 1195.9231
$ws.Range("N138").Value = -23966.207

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 5122.2
$ws.Range("I141").Value = 2226.2727
$ws.Range("J141").Value = 8661.666999999999
$ws.Range("K141").Value = 6678.8181
$ws.Range("L141").Value = 25985.001
$ws.Range("M141").Value = -1498.8181
$ws.Range("N141").Value = -36345.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2908.4
$ws.Range("I32").Value = 2755.9597
$ws.Range("J32").Value = 18000
$ws.Range("K32").Value = 2755.9597
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = -2468.9597
$ws.Range("N32").Value = -18574

$ws.Range("H61").Value = 6341.6865
$ws.Range("I61").Value = 6589.2334
$ws.Range("J61").Value = 4219.857
$ws.Range("K61").Value = 6589.2334
$ws.Range("L61").Value = 4219.857
$ws.Range("M61").Value = -6377.2334
$ws.Range("N61").Value = -4643.857

$ws.Range("H74").Value = 1922.4783
$ws.Range("I74").Value = 1148.3269
$ws.Range("J74").Value = 4290.4707
$ws.Range("K74").Value = 1148.3269
$ws.Range("L74").Value = 4290.4707
$ws.Range("M74").Value = -274.3269
$ws.Range("N74").Value = -6038.4707

$ws.Range("H77").Value = 1922.4783
$ws.Range("I77").Value = 1148.3269
$ws.Range("J77").Value = 4290.4707
$ws.Range("K77").Value = 5741.6345
$ws.Range("L77").Value = 21452.3535
$ws.Range("M77").Value = -1373.6345
$ws.Range("N77").Value = -30188.3535

$ws.Range("H101").Value = 54796.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 54796.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 54796.5
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -61286.5

$ws.Range("H123").Value = 84756.86
$ws.Range("J123").Value = 86366.336
$ws.Range("L123").Value = 86366.336
$ws.Range("N123").Value = -96166.336

$ws.Range("H132").Value = 1991.1605
$ws.Range("I132").Value = 1794.9048
$ws.Range("J132").Value = 2678.0557
$ws.Range("K132").Value = 5384.7144
$ws.Range("L132").Value = 8034.1671
$ws.Range("M132").Value = -2854.7144
$ws.Range("N132").Value = -13094.1671

$ws.Range("H136").Value = 6341.6865
$ws.Range("I136").Value = 6589.2334
$ws.Range("J136").Value = 4219.857
$ws.Range("K136").Value = 19767.7002
$ws.Range("L136").Value = 12659.571
$ws.Range("M136").Value = -17217.7002
$ws.Range("N136").Value = -17759.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 5258.9
$ws.Range("J4").Value = 150
$ws.Range("L4").Value = 150
$ws.Range("N4").Value = -380

$ws.Range("H99").Value = 2979.75
$ws.Range("I99").Value = 2442.9092
$ws.Range("K99").Value = 2442.9092
$ws.Range("M99").Value = -944.9092000000001

$ws.Range("H105").Value = 4882.0454
$ws.Range("I105").Value = 4079.6667
$ws.Range("J105").Value = 5844.9
$ws.Range("K105").Value = 4079.6667
$ws.Range("L105").Value = 5844.9
$ws.Range("M105").Value = -2332.6667
$ws.Range("N105").Value = -9338.9

$ws.Range("H134").Value = 2864.8657
$ws.Range("I134").Value = 2686.1406
$ws.Range("J134").Value = 6677.6665
$ws.Range("K134").Value = 8058.4218
$ws.Range("L134").Value = 20032.9995
$ws.Range("M134").Value = -5523.4218
$ws.Range("N134").Value = -25102.9995

$ws.Range("H140").Value = 166159.86
$ws.Range("J140").Value = 166159.86
$ws.Range("L140").Value = 166159.86
$ws.Range("N140").Value = -176519.86

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H16").Value = 1668.0952
$ws.Range("J16").Value = 2700.8333
$ws.Range("L16").Value = 2700.8333
$ws.Range("N16").Value = -3274.8333

$ws.Range("H58").Value = 2100
$ws.Range("J58").Value = 1935.5714
$ws.Range("L58").Value = 1935.5714
$ws.Range("N58").Value = -2341.5714

$ws.Range("H99").Value = 5772.091
$ws.Range("I99").Value = 5225.16
$ws.Range("K99").Value = 5225.16
$ws.Range("M99").Value = -3727.16

$ws.Range("H113").Value = 1668.0952
$ws.Range("J113").Value = 2700.8333
$ws.Range("L113").Value = 2700.8333
$ws.Range("N113").Value = -7040.8333

$ws.Range("H126").Value = 5772.091
$ws.Range("I126").Value = 5225.16
$ws.Range("K126").Value = 15675.48
$ws.Range("M126").Value = -13205.48

$ws.Range("H132").Value = 1292430.8
$ws.Range("I132").Value = 1668790
$ws.Range("K132").Value = 5006370
$ws.Range("M132").Value = -5003840

$ws.Range("H134").Value = 2056.946
$ws.Range("I134").Value = 989.5357
$ws.Range("J134").Value = 5377.778
$ws.Range("K134").Value = 2968.6071
$ws.Range("L134").Value = 16133.334
$ws.Range("M134").Value = -433.6071000000002
$ws.Range("N134").Value = -21203.334

$ws.Range("H136").Value = 2100
$ws.Range("J136").Value = 1935.5714
$ws.Range("L136").Value = 5806.7142
$ws.Range("N136").Value = -10906.7142

$ws.Range("H141").Value = 731383.6
$ws.Range("J141").Value = 788165.5600000001
$ws.Range("L141").Value = 788165.5600000001
$ws.Range("N141").Value = -798525.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 7062.875
$ws.Range("I43").Value = 500
$ws.Range("K43").Value = 1500
$ws.Range("M43").Value = -1386

$ws.Range("H56").Value = 6952.375
$ws.Range("I56").Value = 6952.375
$ws.Range("K56").Value = 6952.375
$ws.Range("M56").Value = -6422.375

$ws.Range("H122").Value = 1046.4517
$ws.Range("J122").Value = 1421.8182
$ws.Range("L122").Value = 12796.3638
$ws.Range("N122").Value = -17696.3638

$ws.Range("H131").Value = 171184.66
$ws.Range("J131").Value = 1782.2941
$ws.Range("L131").Value = 5346.8823
$ws.Range("N131").Value = -15426.8823

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 72875.25
$ws.Range("I61").Value = 82714.57000000001
$ws.Range("K61").Value = 82714.57000000001
$ws.Range("M61").Value = -82512.57000000001

$ws.Range("H113").Value = 72875.25
$ws.Range("I113").Value = 82714.57000000001
$ws.Range("K113").Value = 82714.57000000001
$ws.Range("M113").Value = -80544.57000000001

$ws.Range("H136").Value = 4573.1567
$ws.Range("I136").Value = 4168.1
$ws.Range("K136").Value = 12504.3
$ws.Range("M136").Value = -9954.300000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 39999
$ws.Range("J31").Value = 39999
$ws.Range("L31").Value = 39999
$ws.Range("N31").Value = -40695

$ws.Range("H57").Value = 101966
$ws.Range("I57").Value = 87999
$ws.Range("K57").Value = 87999
$ws.Range("M57").Value = -87245

$ws.Range("H82").Value = 55091
$ws.Range("I82").Value = 15273
$ws.Range("J82").Value = 75000
$ws.Range("K82").Value = 15273
$ws.Range("L82").Value = 75000
$ws.Range("M82").Value = -14890
$ws.Range("N82").Value = -75766

$ws.Range("H85").Value = 55091
$ws.Range("I85").Value = 15273
$ws.Range("J85").Value = 75000
$ws.Range("K85").Value = 15273
$ws.Range("L85").Value = 75000
$ws.Range("M85").Value = -13947
$ws.Range("N85").Value = -77652

$ws.Range("H96").Value = 3238.625
$ws.Range("I96").Value = 2951.8
$ws.Range("J96").Value = 3716.6667
$ws.Range("K96").Value = 2951.8
$ws.Range("L96").Value = 3716.6667
$ws.Range("M96").Value = -1578.8
$ws.Range("N96").Value = -6462.6667

$ws.Range("H100").Value = 3250
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -5459
$ws.Range("N100").Value = -8082

$ws.Range("H107").Value = 668.5333000000001
$ws.Range("I107").Value = 651.65
$ws.Range("K107").Value = 1954.95
$ws.Range("M107").Value = -34.94999999999982

$ws.Range("H132").Value = 1482.8043
$ws.Range("I132").Value = 1312.2333
$ws.Range("K132").Value = 3936.699900000001
$ws.Range("M132").Value = -1406.699900000001

$ws.Range("H136").Value = 4836.0894
$ws.Range("I136").Value = 5652.9756
$ws.Range("K136").Value = 16958.9268
$ws.Range("M136").Value = -14408.9268
